# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2210
#   *_new  -> *_FV2304
# Also: turn the header row into a real table (Table1) and freeze it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1, columns A:J and L:U) ----------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$oldHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $oldHeaders[$i]
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$newHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $newHeaders[$i]
}

# --- 2. Turn A1:U66 into an actual Excel Table -------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row -------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
